# Atualizacao de bases das ligas (Peru Liga 1): correcao de linhas 175-188
# conforme diff fornecido (re-mapeamento/rotacao de partidas nas linhas 175-177,
# 180-182 e 183-188; linhas 178 e 179 permanecem inalteradas).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 175
$ws.Range("B175").Value = 7302200
$ws.Range("E175").Value = "Carlos Manucci"
$ws.Range("F175").Value = "Deportivo Binacional"
$ws.Range("G175").Value = 3
$ws.Range("I175").Value = 0
$ws.Range("K175").Value = "H"
$ws.Range("L175").Value = 2
$ws.Range("M175").Value = 3.2
$ws.Range("N175").Value = 3.75
$ws.Range("P175").Value = 3.4
$ws.Range("Q175").Value = 4.333
$ws.Range("S175").Value = 1.85
$ws.Range("T175").Value = 1.95
$ws.Range("U175").Value = 2.5
$ws.Range("V175").Value = 1.85
$ws.Range("W175").Value = 1.95
$ws.Range("X175").Value = 0.75
$ws.Range("Z175").Value = -1
$ws.Range("AA175").Value = 0.8500000000000001
$ws.Range("AB175").Value = -1
$ws.Range("AC175").Value = 0.8500000000000001
$ws.Range("AD175").Value = -1

# Row 176
$ws.Range("B176").Value = 7302795
$ws.Range("E176").Value = "Unin Comercio"
$ws.Range("F176").Value = "Deportivo Garcilaso"
$ws.Range("H176").Value = 2
$ws.Range("I176").Value = 1
$ws.Range("K176").Value = "A"
$ws.Range("L176").Value = 2.25
$ws.Range("M176").Value = 3.3
$ws.Range("N176").Value = 2.7
$ws.Range("O176").Value = 1.75
$ws.Range("P176").Value = 3.6
$ws.Range("Q176").Value = 4
$ws.Range("R176").Value = -0.5
$ws.Range("S176").Value = 1.8
$ws.Range("T176").Value = 2
$ws.Range("U176").Value = 2.75
$ws.Range("V176").Value = 1.825
$ws.Range("W176").Value = 1.975
$ws.Range("X176").Value = -1
$ws.Range("Z176").Value = 3
$ws.Range("AB176").Value = 1
$ws.Range("AC176").Value = 0.4125
$ws.Range("AD176").Value = -0.5

# Row 177
$ws.Range("B177").Value = 7302796
$ws.Range("E177").Value = "Sport Huancayo"
$ws.Range("F177").Value = "Sport Boys"
$ws.Range("G177").Value = 1
$ws.Range("H177").Value = 0
$ws.Range("L177").Value = 1.727
$ws.Range("M177").Value = 3.75
$ws.Range("N177").Value = 4.333
$ws.Range("O177").Value = 1.25
$ws.Range("P177").Value = 5.25
$ws.Range("Q177").Value = 10
$ws.Range("R177").Value = -1.75
$ws.Range("S177").Value = 1.925
$ws.Range("T177").Value = 1.875
$ws.Range("U177").Value = 3
$ws.Range("V177").Value = 1.875
$ws.Range("W177").Value = 1.925
$ws.Range("X177").Value = 0.25
$ws.Range("AA177").Value = -1
$ws.Range("AB177").Value = 0.875
$ws.Range("AC177").Value = -1
$ws.Range("AD177").Value = 0.925

# Row 180
$ws.Range("B180").Value = 7384623
$ws.Range("E180").Value = "Sport Boys"
$ws.Range("F180").Value = "Cienciano"
$ws.Range("G180").Value = 2
$ws.Range("O180").Value = 1.833
$ws.Range("P180").Value = 4
$ws.Range("Q180").Value = 3.2
$ws.Range("R180").Value = -0.5
$ws.Range("S180").Value = 1.925
$ws.Range("T180").Value = 1.875
$ws.Range("U180").Value = 3
$ws.Range("V180").Value = 1.925
$ws.Range("W180").Value = 1.875
$ws.Range("X180").Value = 0.833
$ws.Range("AA180").Value = 0.925
$ws.Range("AC180").Value = 0
$ws.Range("AD180").Value = 0

# Row 181
$ws.Range("B181").Value = 7384622
$ws.Range("E181").Value = "Deportivo Municipal"
$ws.Range("F181").Value = "Academia Deportiva Cantolao"
$ws.Range("G181").Value = 1
$ws.Range("H181").Value = 2
$ws.Range("I181").Value = 0
$ws.Range("J181").Value = 0
$ws.Range("K181").Value = "A"
$ws.Range("L181").Value = 1.444
$ws.Range("M181").Value = 4.333
$ws.Range("N181").Value = 7
$ws.Range("O181").Value = 1.5
$ws.Range("P181").Value = 3.75
$ws.Range("Q181").Value = 6
$ws.Range("R181").Value = -1
$ws.Range("S181").Value = 1.825
$ws.Range("T181").Value = 2.025
$ws.Range("U181").Value = 2.75
$ws.Range("V181").Value = 1.875
$ws.Range("W181").Value = 1.975
$ws.Range("X181").Value = -1
$ws.Range("Z181").Value = 5
$ws.Range("AA181").Value = -1
$ws.Range("AB181").Value = 1.025
$ws.Range("AC181").Value = 0.4375
$ws.Range("AD181").Value = -0.5

# Row 182
$ws.Range("B182").Value = 7384624
$ws.Range("E182").Value = "Cesar Vallejo"
$ws.Range("F182").Value = "Cusco FC"
$ws.Range("G182").Value = 3
$ws.Range("H182").Value = 1
$ws.Range("I182").Value = 1
$ws.Range("J182").Value = 1
$ws.Range("K182").Value = "H"
$ws.Range("L182").Value = 2
$ws.Range("M182").Value = 3.4
$ws.Range("N182").Value = 3.5
$ws.Range("O182").Value = 1.45
$ws.Range("P182").Value = 4.2
$ws.Range("Q182").Value = 6.5
$ws.Range("S182").Value = 1.75
$ws.Range("T182").Value = 2.05
$ws.Range("U182").Value = 2.5
$ws.Range("V182").Value = 1.95
$ws.Range("W182").Value = 1.85
$ws.Range("X182").Value = 0.45
$ws.Range("Z182").Value = -1
$ws.Range("AA182").Value = 0.75
$ws.Range("AB182").Value = -1
$ws.Range("AC182").Value = 0.95
$ws.Range("AD182").Value = -1

# Row 183
$ws.Range("B183").Value = 7384630
$ws.Range("E183").Value = "Atletico Grau"
$ws.Range("F183").Value = "Unin Comercio"
$ws.Range("G183").Value = 0
$ws.Range("H183").Value = 1
$ws.Range("I183").Value = 0
$ws.Range("J183").Value = 1
$ws.Range("K183").Value = "A"
$ws.Range("L183").Value = 2.8
$ws.Range("M183").Value = 3.4
$ws.Range("N183").Value = 2.15
$ws.Range("O183").Value = 1.75
$ws.Range("P183").Value = 3.6
$ws.Range("Q183").Value = 3.8
$ws.Range("R183").Value = -0.75
$ws.Range("S183").Value = 2
$ws.Range("T183").Value = 1.8
$ws.Range("U183").Value = 3
$ws.Range("V183").Value = 1.85
$ws.Range("W183").Value = 1.95
$ws.Range("X183").Value = -1
$ws.Range("Z183").Value = 2.8
$ws.Range("AA183").Value = -1
$ws.Range("AB183").Value = 0.8
$ws.Range("AC183").Value = -1
$ws.Range("AD183").Value = 0.95

# Row 184
$ws.Range("B184").Value = 7384627
$ws.Range("E184").Value = "Universitario de Deportes"
$ws.Range("F184").Value = "Sport Huancayo"
$ws.Range("G184").Value = 2
$ws.Range("I184").Value = 1
$ws.Range("K184").Value = "H"
$ws.Range("L184").Value = 1.25
$ws.Range("M184").Value = 5
$ws.Range("N184").Value = 12
$ws.Range("O184").Value = 1.181
$ws.Range("P184").Value = 6
$ws.Range("Q184").Value = 13
$ws.Range("R184").Value = -1.75
$ws.Range("S184").Value = 1.8
$ws.Range("T184").Value = 2
$ws.Range("U184").Value = 2.75
$ws.Range("V184").Value = 1.85
$ws.Range("W184").Value = 1.95
$ws.Range("X184").Value = 0.181
$ws.Range("Y184").Value = -1
$ws.Range("AA184").Value = 0.4
$ws.Range("AB184").Value = -0.5
$ws.Range("AD184").Value = 0.95

# Row 185
$ws.Range("B185").Value = 7384629
$ws.Range("E185").Value = "Deportivo Garcilaso"
$ws.Range("F185").Value = "Alianza Lima"
$ws.Range("G185").Value = 0
$ws.Range("H185").Value = 1
$ws.Range("I185").Value = 0
$ws.Range("J185").Value = 1
$ws.Range("K185").Value = "A"
$ws.Range("L185").Value = 2.625
$ws.Range("M185").Value = 3.3
$ws.Range("N185").Value = 2.5
$ws.Range("O185").Value = 2.7
$ws.Range("P185").Value = 3.4
$ws.Range("Q185").Value = 2.375
$ws.Range("R185").Value = 0
$ws.Range("S185").Value = 2.025
$ws.Range("T185").Value = 1.775
$ws.Range("U185").Value = 2.25
$ws.Range("V185").Value = 1.825
$ws.Range("W185").Value = 1.975
$ws.Range("X185").Value = -1
$ws.Range("Z185").Value = 1.375
$ws.Range("AA185").Value = -1
$ws.Range("AB185").Value = 0.7749999999999999
$ws.Range("AD185").Value = 0.9750000000000001

# Row 186
$ws.Range("B186").Value = 7384626
$ws.Range("E186").Value = "Sporting Cristal"
$ws.Range("F186").Value = "Alianza Atletico"
$ws.Range("G186").Value = 3
$ws.Range("H186").Value = 0
$ws.Range("I186").Value = 3
$ws.Range("J186").Value = 0
$ws.Range("K186").Value = "H"
$ws.Range("L186").Value = 1.3
$ws.Range("M186").Value = 5
$ws.Range("N186").Value = 9
$ws.Range("O186").Value = 1.166
$ws.Range("P186").Value = 6.5
$ws.Range("Q186").Value = 13
$ws.Range("R186").Value = -2
$ws.Range("S186").Value = 1.85
$ws.Range("T186").Value = 1.95
$ws.Range("U186").Value = 3.25
$ws.Range("V186").Value = 2
$ws.Range("W186").Value = 1.8
$ws.Range("X186").Value = 0.1659999999999999
$ws.Range("Z186").Value = -1
$ws.Range("AA186").Value = 0.8500000000000001
$ws.Range("AB186").Value = -1
$ws.Range("AC186").Value = -0.5
$ws.Range("AD186").Value = 0.4

# Row 187
$ws.Range("B187").Value = 7384625
$ws.Range("E187").Value = "AD Tarma"
$ws.Range("F187").Value = "Carlos Manucci"
$ws.Range("H187").Value = 0
$ws.Range("J187").Value = 0
$ws.Range("K187").Value = "D"
$ws.Range("L187").Value = 1.5
$ws.Range("M187").Value = 3.75
$ws.Range("N187").Value = 7
$ws.Range("O187").Value = 1.363
$ws.Range("P187").Value = 4.333
$ws.Range("Q187").Value = 9.5
$ws.Range("R187").Value = -1.25
$ws.Range("S187").Value = 1.875
$ws.Range("T187").Value = 1.925
$ws.Range("U187").Value = 2.5
$ws.Range("V187").Value = 1.8
$ws.Range("W187").Value = 2
$ws.Range("Y187").Value = 3.333
$ws.Range("Z187").Value = -1
$ws.Range("AB187").Value = 0.925
$ws.Range("AD187").Value = 1

# Row 188
$ws.Range("B188").Value = 7384628
$ws.Range("E188").Value = "Deportivo Binacional"
$ws.Range("F188").Value = "FBC Melgar"
$ws.Range("G188").Value = 1
$ws.Range("H188").Value = 2
$ws.Range("I188").Value = 1
$ws.Range("L188").Value = 2.75
$ws.Range("N188").Value = 2.375
$ws.Range("O188").Value = 3.3
$ws.Range("P188").Value = 3.6
$ws.Range("Q188").Value = 2
$ws.Range("R188").Value = 0.5
$ws.Range("S188").Value = 1.8
$ws.Range("T188").Value = 2
$ws.Range("U188").Value = 2.75
$ws.Range("V188").Value = 1.975
$ws.Range("W188").Value = 1.875
$ws.Range("Z188").Value = 1
$ws.Range("AB188").Value = 1
$ws.Range("AC188").Value = 0.4875
$ws.Range("AD188").Value = -0.5
